$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-10-03 -> 2023-10-04, serials 45202 -> 45203) for every data row
# (rows 2 through 120).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 120) { $lastRow = 120 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
